$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.047.61"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.828.68"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Formula = "'0.9997"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Formula = "'243.10"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Formula = "'0.6226"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").Formula = "'1.002"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Formula = "'0.07430"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").Formula = "'0.2915"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Formula = "'23.09"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").Formula = "'0.07696"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "1.832.27"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Formula = "'4.996"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Formula = "'0.6675"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Formula = "'82.36"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Formula = "'0.000009271"
$ws.Range("E16").Value = "  -7.65%  "
$ws.Range("D17").Formula = "'5.915"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").Value = "29.066.69"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "2.076.94"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Formula = "'12.58"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Formula = "'219.72"
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("D22").Formula = "'1.003"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Formula = "'7.098"
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").Formula = "'1.003"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Formula = "'160.14"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Formula = "'0.1387"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Formula = "'8.482"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Formula = "'17.82"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Formula = "'1.492"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Formula = "'0.05662"
$ws.Range("E30").Value = "  +7.59%  "
$ws.Range("D31").Formula = "'4.160"
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").Formula = "'4.106"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").Formula = "'1.209"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Formula = "'0.7383"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Formula = "'1.819"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").Formula = "'1.137"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Formula = "'2.677"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Formula = "'2.765"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "1.218.53"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("D40").Formula = "'0.01766"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Formula = "'6.466"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Formula = "'0.8897"
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("D43").Formula = "'1.003"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Formula = "'101.69"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "1.977.89"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Formula = "'66.14"
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Formula = "'0.00000000123"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Formula = "'0.5086"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Formula = "'0.07464"
$ws.Range("E49").Value = "  +16.92%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Formula = "'8.991"
$ws.Range("E51").Value = "  +0.48%  "
